# Update table Sp3: training set composition of other published ML algorithms
# The data rows (3-16) are reordered: the "2017_Tsai" (CIRCLE-seq) and
# "2018_Cameron" (SITE-seq) rows move up from the bottom of the short block
# to right after "2015_Slaymaker", and the merged "2015_Haeussler" block
# (7 rows) shifts down from A8:A14 to A10:A16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the old merged range first so the member cells can be edited
# independently, then we will re-merge the new range once data is in place.
$ws.Range("A8:A14").UnMerge()

# --- Row 3: now 2015_Slaymaker / BLESS / In-vitro DSB ---
$ws.Range("B3").Value = '2015_Slaymaker'
$ws.Range("C3").Value = 'BLESS'
$ws.Range("D3").Value = 'In-vitro DSB'
$ws.Range("E3").ClearContents()
$ws.Range("G3").Value = 'Y'
$ws.Range("I3").ClearContents()

# --- Row 4: now 2017_Tsai / CIRCLE-seq ---
$ws.Range("B4").Value = '2017_Tsai'
$ws.Range("C4").Value = 'CIRCLE-seq'
$ws.Range("F4").Value = 'Y'
$ws.Range("G4").ClearContents()
$ws.Range("I4").Value = 'Y'

# --- Row 5: now 2018_Cameron / SITE-seq ---
$ws.Range("B5").Value = '2018_Cameron'
$ws.Range("C5").Value = 'SITE-seq'
$ws.Range("D5").Value = 'In-vitro DSB'
$ws.Range("E5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("I5").Value = 'Y'

# --- Row 6: now 2015_Tsai / GUIDE-seq ---
$ws.Range("B6").Value = '2015_Tsai'
$ws.Range("C6").Value = 'GUIDE-seq'
$ws.Range("D6").Value = 'in-vivo tag intergration'
$ws.Range("F6").Value = 'Y'
$ws.Range("G6").Value = 'Y'
$ws.Range("H6").Value = 'Independent testing set'

# --- Row 7: now 2016_Kleinstiver / GUIDE-seq ---
$ws.Range("B7").Value = '2016_Kleinstiver'
$ws.Range("E7").Value = 'Independent testing set'
$ws.Range("F7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = 'Independent testing set'

# --- Row 8: now 2018_Listgarten / GUIDE-seq (the "header-like" row) ---
$ws.Range("A8").ClearContents()
$ws.Range("B8").Value = '2018_Listgarten'
$ws.Range("C8").Value = 'GUIDE-seq'
$ws.Range("D8").Value = 'in-vivo tag intergration'
$ws.Range("E8").Value = 'Independent testing set'
$ws.Range("F8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 'Independent testing set'

# --- Row 9: now 2016_Doench / CFD ---
$ws.Range("B9").Value = '2016_Doench'
$ws.Range("C9").Value = 'CFD'
$ws.Range("D9").Value = 'In-vivo gene knockout'
$ws.Range("E9").Value = 'Y'
$ws.Range("F9").ClearContents()
$ws.Range("H9").ClearContents()

# --- Row 10 (was 8): now the start of the Haeussler merged block, 2013_Hsu ---
$ws.Range("A10").Value = '2015_Haeussler (collection of 7 studies)'
$ws.Range("B10").Value = '2013_Hsu'
$ws.Range("C10").Value = 'Targeted sequencing '
$ws.Range("D10").Value = 'in-vivo indels'
$ws.Range("G10").ClearContents()

# --- Row 11 (was 9): 2014_Cho ---
$ws.Range("B11").Value = '2014_Cho'
$ws.Range("C11").Value = 'Targeted sequencing '
$ws.Range("D11").Value = 'in-vivo indels'

# --- Row 12 (was 10): 2015_Frock ---
$ws.Range("B12").Value = '2015_Frock'
$ws.Range("C12").Value = 'HTGTS'
$ws.Range("G12").Value = 'Y'

# --- Row 13 (was 11): 2015_Kim ---
$ws.Range("B13").Value = '2015_Kim'
$ws.Range("C13").Value = 'DIGENOME-seq'
$ws.Range("G13").ClearContents()

# --- Row 14 (was 12): 2015_Wang ---
$ws.Range("B14").Value = '2015_Wang'
$ws.Range("C14").Value = 'Lentiviral Integration + Targeted sequencing'
$ws.Range("D14").Value = 'in-vivo indels'

# --- Row 15 (was 13): 2015_Ran ---
$ws.Range("B15").Value = '2015_Ran'
$ws.Range("C15").Value = 'BLESS'
$ws.Range("G15").Value = 'Y'
$ws.Range("H15").Value = 'Y(GUIDE-seq data excluded)'

# --- Row 16 (was 14): 2016_Kim ---
$ws.Range("B16").Value = '2016_Kim'
$ws.Range("C16").Value = 'DIGENOME-seq2'
$ws.Range("F16").Value = 'Y'
$ws.Range("H16").Value = 'Y(GUIDE-seq data excluded)'

# Re-merge column A for the Haeussler block at its new location
$ws.Range("A10:A16").Merge()

# Update the selected cell to match the saved selection in the edited file
$ws.Range("D12").Select()
